# Applies the "all factures in one xml file" edit to the EK RCH sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EK RCH")

# A15 previously held a real date serial (44351); it is replaced with the
# literal text "03.06.2021" (matching the text already used in A16 and
# other rows in this column). Copying the already-correct A16 cell keeps
# the cell's existing style/number-format intact instead of Excel's
# auto-detect-as-date behaviour that a plain Value assignment would cause.
$ws.Range("A16").Copy($ws.Range("A15"))

# New column P ("pořadové číslo" / running row number) is populated for
# rows 15-69 with sequential integers 1..55. Row 16 already had a
# placeholder value (123456789) that gets replaced by the correct
# sequential number (2).
$row = 15
$n = 1
while ($row -le 69) {
    $ws.Cells.Item($row, 16).Value = $n
    $row++
    $n++
}
